# Append the 2026-02-07 price-scrape row (row 29) to the tracking sheet.
# Values are written as text (matching the existing rows, which all store
# their Date/Price/Discount/Incredible entries as shared strings rather
# than numbers/dates) by entering them as formulas that evaluate to text
# and then converting those formulas to static values via copy/paste
# special. This avoids Excel's automatic number/date type-inference (and
# the "quote prefix" cell style it would otherwise stamp onto the cells)
# that a direct .Value = "2026-02-07" assignment would trigger.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A29").Formula = '="2026-02-07"'
$ws.Range("B29").Formula = '="21700000"'
$ws.Range("C29").Formula = '="0"'
$ws.Range("D29").Formula = '="0"'

$ws.Range("A29:D29").Copy()
$ws.Range("A29:D29").PasteSpecial(-4163)
